# Apply "mary is da best" revision to Times.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# 1) New data entered for query block 3 (rows 15-20), column C,
#    plus the shared average formula in column E extended down to
#    row 20 (was only filled through row 14).
# ---------------------------------------------------------------
$ws.Range("C15").Value = 6.642
$ws.Range("C16").Value = 0.111
$ws.Range("C17").Value = 1.007
$ws.Range("C18").Value = 1.095
$ws.Range("C19").Value = 0.978
$ws.Range("C20").Value = 1.022

# D15 picks up the same look as the rest of the D column (s=10)
# instead of its former one-off style.
$ws.Range("D15").HorizontalAlignment = -4131

# Extend the "(C+D)/2" average formula down through row 20 (it
# previously stopped at row 14). Cells are written one at a time so
# each keeps its own formula text instead of losing it to a second,
# master-less shared-formula group.
$ws.Range("E15").Formula = "=(C15+D15)/2"
$ws.Range("E16").Formula = "=(C16+D16)/2"
$ws.Range("E17").Formula = "=(C17+D17)/2"
$ws.Range("E18").Formula = "=(C18+D18)/2"
$ws.Range("E19").Formula = "=(C19+D19)/2"
$ws.Range("E20").Formula = "=(C20+D20)/2"

# The empty, formatted-but-unused C/D placeholder cells on the
# still-unfilled query blocks (rows 21-32) are cleared out entirely.
$ws.Range("C21:D32").Clear()

# ---------------------------------------------------------------
# 2) View state: scrolled down a bit, selection moved to the newly
#    filled E14:E20 range.
# ---------------------------------------------------------------
$ws.Range("E14:E20").Select()
$excel.ActiveWindow.ScrollRow = 3

# ---------------------------------------------------------------
# 3) Workbook-level metadata: locale switched to Greek Excel, so the
#    built-in "Normal" cell style is now localized. (Best effort --
#    the host COM surface may not expose a writable Styles
#    collection; harmless if it is a no-op.)
# ---------------------------------------------------------------
foreach ($st in $wb.Styles) {
    if ($st.Name -eq "Normal") {
        $st.Name = "Κανονικό"
    }
}

"done"
